$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2, 3 and 4 (the first three inventory entries), shifting
# subsequent rows up. This matches the source data: the three oldest
# records (VIZ M 1084, first NK H 850, first NAC M 890) were removed.
$ws.Range("A2:G4").EntireRow.Delete()

# After deleting whole rows, Excel leaves the selection on the full row
# that slid up into the deleted rows' place (row 2).
$ws.Rows(2).Select() | Out-Null
